{"js": "// 1. Update the date heading (first paragraph of the body).\nconst dateResults = context.document.body.search(\"2025-11-25 Tuesday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"2025-11-26 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2. Update the answers table: same 20-row x 5-column shape, only the text of\n// the previously populated cells changes (new multiplication problems/answers).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Row 1 (index 0)\ntable.getCell(0, 0).value = \"729\u00d75=3645\";\ntable.getCell(0, 1).value = \"868\u00d78=6944\";\ntable.getCell(0, 2).value = \"244\u00d74=976\";\ntable.getCell(0, 3).value = \"332\u00d72=664\";\ntable.getCell(0, 4).value = \"732\u00d78=5856\";\n\n// Row 5 (index 4)\ntable.getCell(4, 0).value = \"738\u00d76=4428\";\ntable.getCell(4, 1).value = \"187\u00d76=1122\";\ntable.getCell(4, 2).value = \"867\u00d74=3468\";\ntable.getCell(4, 3).value = \"649\u00d76=3894\";\ntable.getCell(4, 4).value = \"543\u00d77=3801\";\n\n// Row 10 (index 9)\ntable.getCell(9, 0).value = \"904\u00d75=4520\";\ntable.getCell(9, 1).value = \"441\u00d73=1323\";\ntable.getCell(9, 2).value = \"353\u00d72=706\";\ntable.getCell(9, 3).value = \"346\u00d79=3114\";\ntable.getCell(9, 4).value = \"849\u00d77=5943\";\n\n// Row 15 (index 14)\ntable.getCell(14, 0).value = \"544\u00d75=2720\";\ntable.getCell(14, 1).value = \"313\u00d76=1878\";\ntable.getCell(14, 2).value = \"863\u00d74=3452\";\ntable.getCell(14, 3).value = \"958\u00d75=4790\";\ntable.getCell(14, 4).value = \"563\u00d76=3378\";\n\n// Row 20 (index 19)\ntable.getCell(19, 0).value = \"671\u00d75=3355\";\ntable.getCell(19, 1).value = \"348\u00d78=2784\";\ntable.getCell(19, 2).value = \"719\u00d78=5752\";\ntable.getCell(19, 3).value = \"150\u00d76=900\";\ntable.getCell(19, 4).value = \"414\u00d79=3726\";\n\nawait context.sync();\n", "ps1": "# 1. Update the date heading (first paragraph of the body).\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.Text = \"2025-11-25 Tuesday\"\n$find.Replacement.Text = \"2025-11-26 Wednesday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2. Update the answers table: same 20-row x 5-column shape, only the text of\n# the previously populated cells changes (new multiplication problems/answers).\n$table = $d.Tables.Item(1)\n\n# Row 1\n$table.Cell(1, 1).Range.Text = \"729\u00d75=3645\"\n$table.Cell(1, 2).Range.Text = \"868\u00d78=6944\"\n$table.Cell(1, 3).Range.Text = \"244\u00d74=976\"\n$table.Cell(1, 4).Range.Text = \"332\u00d72=664\"\n$table.Cell(1, 5).Range.Text = \"732\u00d78=5856\"\n\n# Row 5\n$table.Cell(5, 1).Range.Text = \"738\u00d76=4428\"\n$table.Cell(5, 2).Range.Text = \"187\u00d76=1122\"\n$table.Cell(5, 3).Range.Text = \"867\u00d74=3468\"\n$table.Cell(5, 4).Range.Text = \"649\u00d76=3894\"\n$table.Cell(5, 5).Range.Text = \"543\u00d77=3801\"\n\n# Row 10\n$table.Cell(10, 1).Range.Text = \"904\u00d75=4520\"\n$table.Cell(10, 2).Range.Text = \"441\u00d73=1323\"\n$table.Cell(10, 3).Range.Text = \"353\u00d72=706\"\n$table.Cell(10, 4).Range.Text = \"346\u00d79=3114\"\n$table.Cell(10, 5).Range.Text = \"849\u00d77=5943\"\n\n# Row 15\n$table.Cell(15, 1).Range.Text = \"544\u00d75=2720\"\n$table.Cell(15, 2).Range.Text = \"313\u00d76=1878\"\n$table.Cell(15, 3).Range.Text = \"863\u00d74=3452\"\n$table.Cell(15, 4).Range.Text = \"958\u00d75=4790\"\n$table.Cell(15, 5).Range.Text = \"563\u00d76=3378\"\n\n# Row 20\n$table.Cell(20, 1).Range.Text = \"671\u00d75=3355\"\n$table.Cell(20, 2).Range.Text = \"348\u00d78=2784\"\n$table.Cell(20, 3).Range.Text = \"719\u00d78=5752\"\n$table.Cell(20, 4).Range.Text = \"150\u00d76=900\"\n$table.Cell(20, 5).Range.Text = \"414\u00d79=3726\"\n"}
